$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data update -----------------------------------------------------------
# The "purchase" module's route was renamed to live in "purchase_stock"
# (Odoo module split). Update every route_ids (column R) cell that still
# references the old "purchase.route_warehouse0_buy" xmlid so it points at
# "purchase_stock.route_warehouse0_buy" instead - whether it appears alone
# or combined with other routes in a comma separated list.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 18).End(-4162).Row   # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 30 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 18)   # column R = route_ids
    $val = $cell.Value()
    if ($val -ne $null -and $val.ToString().Contains("purchase.route_warehouse0_buy")) {
        $cell.Value = $val.ToString().Replace("purchase.route_warehouse0_buy", "purchase_stock.route_warehouse0_buy")
    }
}

# --- View / selection update -------------------------------------------------
# Move the frozen-pane viewport: the sheet stays split after column C / row 1,
# but the window had been scrolled so the bottom-right pane's visible corner
# was R17; re-anchor it at R2 (i.e. scroll back up) and update the active
# selections in each pane accordingly.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("R2").Select() | Out-Null
$win.FreezePanes = $true

$ws.Range("A2").Select() | Out-Null
$ws.Range("R2").Select() | Out-Null
